$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "fecha_venta"
$ws.Range("E2:E3").NumberFormat = "@"
$ws.Range("E2").Value = "2022-03-21"
$ws.Range("E3").Value = "2022-03-21"
